$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new value in B4 (stored as a shared string "नमस्ते")
$ws.Range("B4").Value = "नमस्ते"

# Match the centered style used by the rest of column B/the data block
$ws.Range("B4").HorizontalAlignment = $ws.Range("B3").HorizontalAlignment

# Select B4 so it becomes the active cell / selection on the sheet
$ws.Activate()
$ws.Range("B4").Select()

# Add a new defined name with a non-ASCII (Greek) name pointing at Sheet1!$B$4.
# The COM method resolver for Names.Add needs the Name argument to start with
# an ASCII character, so add it with a placeholder ASCII name first and then
# rename it to the desired Unicode name.
$newName = $wb.Names.Add("TempName_ForRename", "=Sheet1!`$B`$4")
$newName.Name = "Χαιρετισμός"
